$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.02907897629796788
$ws.Range("C2").Value = 0.3131278957257717
$ws.Range("D2").Value = 0.181524606355785
$ws.Range("E2").Value = 0.4260570458938391
$ws.Range("F2").Value = 0.43998257208981
$ws.Range("G2").Value = 15

$ws.Range("B3").Value = 0.0370360702193368
$ws.Range("C3").Value = 0.2877205798324301
$ws.Range("D3").Value = 0.135345957066826
$ws.Range("E3").Value = 0.3678939481247633
$ws.Range("F3").Value = 0.3798420736582994
$ws.Range("G3").Value = 14

$ws.Range("B4").Value = 0.01989109342689383
$ws.Range("C4").Value = 0.3298745939743749
$ws.Range("D4").Value = 0.1503185731831121
$ws.Range("E4").Value = 0.3877093926939507
$ws.Range("F4").Value = 0.4030092977225764
$ws.Range("G4").Value = 13

$ws.Range("B5").Value = 0.06273519545459039
$ws.Range("C5").Value = 0.2460334154802962
$ws.Range("D5").Value = 0.09558537405740077
$ws.Range("E5").Value = 0.3091688439306276
$ws.Range("F5").Value = 0.3161984462764572
$ws.Range("G5").Value = 12

$ws.Range("B6").Value = 0.02711121047383215
$ws.Range("C6").Value = 0.2779584568479502
$ws.Range("D6").Value = 0.1216586105046287
$ws.Range("E6").Value = 0.34879594393374
$ws.Range("F6").Value = 0.3647135205176791
$ws.Range("G6").Value = 11

$ws.Range("B7").Value = -0.01088135635153479
$ws.Range("C7").Value = 0.2849421119723689
$ws.Range("D7").Value = 0.09825224066429059
$ws.Range("E7").Value = 0.3134521345664926
$ws.Range("F7").Value = 0.3302084135617004
$ws.Range("G7").Value = 10

$ws.Range("B8").Value = -0.04428645741563344
$ws.Range("C8").Value = 0.3765628100937468
$ws.Range("D8").Value = 0.1690784592812682
$ws.Range("E8").Value = 0.4111915116843588
$ws.Range("F8").Value = 0.4335975266214011
$ws.Range("G8").Value = 9

$ws.Range("B9").Value = -0.01393931246739222
$ws.Range("C9").Value = 0.3509192590318558
$ws.Range("D9").Value = 0.1607817096315818
$ws.Range("E9").Value = 0.4009759464501354
$ws.Range("F9").Value = 0.4284022211487281
$ws.Range("G9").Value = 8

$ws.Range("B10").Value = -0.03688841855209302
$ws.Range("C10").Value = 0.2799590153621541
$ws.Range("D10").Value = 0.1329006694222477
$ws.Range("E10").Value = 0.3645554408073589
$ws.Range("F10").Value = 0.3917438359423185
$ws.Range("G10").Value = 7

$ws.Range("B11").Value = -0.06364682135181432
$ws.Range("C11").Value = 0.1670219060428917
$ws.Range("D11").Value = 0.03807752624502202
$ws.Range("E11").Value = 0.1951346362002964
$ws.Range("F11").Value = 0.2020691219662191
$ws.Range("G11").Value = 6
